# Add new columns (K, L) of test data to the "addReseller" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("addReseller")

# Fill column K (header + data), then column L, so new shared strings are
# created in the same order they appear in the saved workbook.
$ws.Range("K1").Value = "userIDMPassword"
$ws.Range("K2").Value = "Independent12#"
$ws.Range("L1").Value = "resNameSuf"
$ws.Range("L2").Value = "& brother's"

# Match the yellow header fill used by the rest of row 1.
$ws.Range("K1").Interior.Color = 65535
$ws.Range("L1").Interior.Color = 65535

# Update the selected cell on this sheet, as recorded after the edit.
$ws.Activate()
$ws.Range("E8").Select()
